# AT7215KV200.xlsx update
#  - Renames propeller labels from "NN*MM" notation to "NNxMM" notation
#    (e.g. "19*10" -> "19x10") and reorders "Propeller" to the front of the
#    shared-string table (matches how Excel re-levels the shared strings
#    once the header cell C1 is re-saved referencing it first).
#  - Updates the window selection left over from editing (now on D79,
#    having scrolled the view down near row 55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - "Propeller" (text unchanged, just touched so the shared
# string table keeps it as the first-seen unique string, matching the
# target workbook's string order).
$ws.Range("C1").Value = "Propeller"

# Row groups for each propeller, with "*" swapped for "x".
for ($r = 2;  $r -le 12; $r++) { $ws.Range("C$r").Value = "19x10"  }
for ($r = 13; $r -le 23; $r++) { $ws.Range("C$r").Value = "19x10b" }
for ($r = 24; $r -le 34; $r++) { $ws.Range("C$r").Value = "20x10"  }
for ($r = 35; $r -le 45; $r++) { $ws.Range("C$r").Value = "20x10b" }
for ($r = 46; $r -le 56; $r++) { $ws.Range("C$r").Value = "20x10c" }
for ($r = 57; $r -le 67; $r++) { $ws.Range("C$r").Value = "20x10d" }
for ($r = 68; $r -le 78; $r++) { $ws.Range("C$r").Value = "21x10"  }
for ($r = 79; $r -le 89; $r++) { $ws.Range("C$r").Value = "21x10b" }

# Restore the saved window selection/scroll position.
$ws.Range("D79").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
